# Apply the data updates described by the commit ("Add files via upload")
# to the Albuwell PFAS sample sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a data-entry typo: C16 should read "Ladder" (matching rows 14-15),
# not "Sample".
$ws.Range("C16").Value = "Ladder"

# --- Fill in / correct the H1:H4 (rows 86-89) dilution block with the
# "4 D6" / "8 D6" sample labels, updated dilution (20) and replicate
# numbers.
$ws.Range("B86").Value = "4 D6"
$ws.Range("C86").Value = "Sample"
$ws.Range("E86").Value = 20
$ws.Range("F86").Value = 1

$ws.Range("B87").Value = "4 D6"
$ws.Range("C87").Value = "Sample"
$ws.Range("E87").Value = 20
$ws.Range("F87").Value = 2

$ws.Range("B88").Value = "8 D6"
$ws.Range("C88").Value = "Sample"
$ws.Range("E88").Value = 20
$ws.Range("F88").Value = 1

$ws.Range("B89").Value = "8 D6"
$ws.Range("C89").Value = "Sample"
$ws.Range("E89").Value = 20
$ws.Range("F89").Value = 2

# --- Rows H6 and H7 (rows 91-92) get fully cleared (contents + formatting)
[void]$ws.Range("B91:F91").Clear()
[void]$ws.Range("B92:F92").Clear()

# --- Update the on-screen selection to match where the author left off.
[void]$ws.Range("A1:F97").Select()
[void]$ws.Range("C16").Activate()
